# Structural analysis example work: trim stray leading/trailing
# whitespace out of the label cells on both sheets (the shared-string
# table collapses the now-duplicate entries automatically on save),
# restore the selection to the range the author was working in, and
# renumber the first printed page of the "connectiosn" sheet.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $val = $cell.Value2
        if ($val -ne $null -and $val.GetType().Name -eq "String") {
            $trimmed = $val.Trim()
            if ($trimmed -ne $val) {
                $cell.Value = $trimmed
            }
        }
    }
}

# Restore the multi-area selection that was active while doing this
# clean-up pass on both sheets.
$wsStruct = $wb.Worksheets.Item("struct")
$wsStruct.Range("A2:A33,G1,L1:AG1").Select()

$wsConn = $wb.Worksheets.Item("connectiosn")
$wsConn.Activate()
$wsConn.Range("A2:A33,G1,L1:AG1").Select()

# The "connectiosn" sheet's printed output should start counting pages
# from 0 (previously defaulted to 1).
$wsConn.PageSetup.FirstPageNumber = 0
